# Auto-generated Excel COM-interop script
# Applies the Cerberus_Profits market-data refresh described in the commit diff.
# For each changed row the H/I/J/K/L/M/N ("currentAveragePrice*", "LevePrice*", "LeveProfit*")
# columns are refreshed with newly-fetched market values. A few rows also gain or lose
# an M/N cell entirely (sparse LeveProfit columns), handled via ClearContents()/Value writes.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1267.7
$ws.Range("I32").Value = 1250.5
$ws.Range("J32").Value = 1272
$ws.Range("K32").Value = 1250.5
$ws.Range("L32").Value = 1272
$ws.Range("M32").Value = -924.5
$ws.Range("N32").Value = -1924
# Row 40
$ws.Range("H40").Value = 1787.3334
$ws.Range("I40").Value = 1056.1666
$ws.Range("J40").Value = 3249.6667
$ws.Range("K40").Value = 1056.1666
$ws.Range("L40").Value = 3249.6667
$ws.Range("M40").Value = -881.1666
$ws.Range("N40").Value = -3599.6667
# Row 55
$ws.Range("H55").Value = 97.59090999999999
$ws.Range("I55").Value = 95
$ws.Range("J55").Value = 97.84999999999999
$ws.Range("K55").Value = 95
$ws.Range("L55").Value = 97.84999999999999
$ws.Range("M55").Value = 119
$ws.Range("N55").Value = -525.85
# Row 113
$ws.Range("H113").Value = 6718.9443
$ws.Range("I113").Value = 6375.4
$ws.Range("J113").Value = 6851.077
$ws.Range("K113").Value = 6375.4
$ws.Range("L113").Value = 6851.077
$ws.Range("M113").Value = -3121.4
$ws.Range("N113").Value = -13359.077
# Row 116
$ws.Range("H116").Value = 9016.107
$ws.Range("I116").Value = 10053.615
$ws.Range("J116").Value = 8116.933
$ws.Range("K116").Value = 10053.615
$ws.Range("L116").Value = 8116.933
$ws.Range("M116").Value = -6611.615
$ws.Range("N116").Value = -15000.933
# Row 132
$ws.Range("H132").Value = 2905.7666
$ws.Range("I132").Value = 2468.2307
$ws.Range("J132").Value = 5749.75
$ws.Range("K132").Value = 7404.6921
$ws.Range("L132").Value = 17249.25
$ws.Range("M132").Value = -4874.6921
$ws.Range("N132").Value = -22309.25
# Row 137
$ws.Range("H137").Value = 2589.276
$ws.Range("I137").Value = 1956.5714
$ws.Range("J137").Value = 3179.8
$ws.Range("K137").Value = 5869.7142
$ws.Range("L137").Value = 9539.400000000001
$ws.Range("M137").Value = -3319.7142
$ws.Range("N137").Value = -14639.4
# Row 138
$ws.Range("H138").Value = 9356
$ws.Range("I138").Value = 8239.4
$ws.Range("J138").Value = 9863.546
$ws.Range("K138").Value = 24718.2
$ws.Range("L138").Value = 29590.638
$ws.Range("M138").Value = -19578.2
$ws.Range("N138").Value = -39870.638

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1763.65
$ws.Range("I2").Value = 1394.2667
$ws.Range("K2").Value = 1394.2667
$ws.Range("M2").Value = -1281.2667
# Row 23
$ws.Range("H23").Value = 59000
$ws.Range("J23").Value = 59000
$ws.Range("L23").Value = 59000
$ws.Range("N23").Value = -59518
# Row 32
$ws.Range("H32").Value = 1309.1578
$ws.Range("I32").Value = 882.1539
$ws.Range("K32").Value = 882.1539
$ws.Range("M32").Value = -595.1539
# Row 44
$ws.Range("H44").Value = 47500
$ws.Range("J44").Value = 47500
$ws.Range("L44").Value = 47500
$ws.Range("N44").Value = -48476
# Row 45
$ws.Range("H45").Value = 1830.6875
$ws.Range("I45").Value = 946.63635
$ws.Range("J45").Value = 3775.6
$ws.Range("K45").Value = 946.63635
$ws.Range("L45").Value = 3775.6
$ws.Range("M45").Value = -569.63635
$ws.Range("N45").Value = -4529.6
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
# Row 61
$ws.Range("H61").Value = 6314.516
$ws.Range("I61").Value = 4541.6665
$ws.Range("J61").Value = 12392.857
$ws.Range("K61").Value = 4541.6665
$ws.Range("L61").Value = 12392.857
$ws.Range("M61").Value = -4329.6665
$ws.Range("N61").Value = -12816.857
# Row 113
$ws.Range("H113").Value = 99999
$ws.Range("J113").Value = 99999
$ws.Range("L113").Value = 99999
$ws.Range("N113").Value = -108677
# Row 116
$ws.Range("H116").Value = 1763.65
$ws.Range("I116").Value = 1394.2667
$ws.Range("K116").Value = 1394.2667
$ws.Range("M116").Value = 899.7333000000001
# Row 136
$ws.Range("H136").Value = 6314.516
$ws.Range("I136").Value = 4541.6665
$ws.Range("J136").Value = 12392.857
$ws.Range("K136").Value = 13624.9995
$ws.Range("L136").Value = 37178.571
$ws.Range("M136").Value = -11074.9995
$ws.Range("N136").Value = -42278.571

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1763.65
$ws.Range("I3").Value = 1394.2667
$ws.Range("K3").Value = 1394.2667
$ws.Range("M3").Value = -1280.2667
# Row 86
$ws.Range("H86").Value = 5638.222
$ws.Range("I86").Value = 5054.4
$ws.Range("J86").Value = 6368
$ws.Range("K86").Value = 5054.4
$ws.Range("L86").Value = 6368
$ws.Range("M86").Value = -3931.4
$ws.Range("N86").Value = -8614
# Row 89
$ws.Range("H89").Value = 5638.222
$ws.Range("I89").Value = 5054.4
$ws.Range("J89").Value = 6368
$ws.Range("K89").Value = 25272
$ws.Range("L89").Value = 31840
$ws.Range("M89").Value = -19656
$ws.Range("N89").Value = -43072

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 464999.5
$ws.Range("J9").Value = 464999.5
$ws.Range("L9").Value = 464999.5
$ws.Range("N9").Value = -465335.5
# Row 22
$ws.Range("H22").Value = 1015.7778
$ws.Range("I22").Value = 603.5
$ws.Range("J22").Value = 1345.6
$ws.Range("K22").Value = 603.5
$ws.Range("L22").Value = 1345.6
$ws.Range("M22").Value = -253.5
$ws.Range("N22").Value = -2045.6
# Row 31
$ws.Range("H31").Value = 3967.3958
$ws.Range("I31").Value = 3893
$ws.Range("J31").Value = 3982.275
$ws.Range("K31").Value = 3893
$ws.Range("L31").Value = 3982.275
$ws.Range("M31").Value = -3598
$ws.Range("N31").Value = -4572.275
# Row 34
$ws.Range("H34").Value = 3967.3958
$ws.Range("I34").Value = 3893
$ws.Range("J34").Value = 3982.275
$ws.Range("K34").Value = 3893
$ws.Range("L34").Value = 3982.275
$ws.Range("M34").Value = -3691
$ws.Range("N34").Value = -4386.275
# Row 132
$ws.Range("H132").Value = 2718.5454
$ws.Range("I132").Value = 2609.9048
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7829.714399999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5299.714399999999
$ws.Range("N132").Value = -20060

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3400.6667
$ws.Range("J68").Value = 3574.4211
$ws.Range("L68").Value = 10723.2633
$ws.Range("N68").Value = -12345.2633
# Row 71
$ws.Range("H71").Value = 3400.6667
$ws.Range("J71").Value = 3574.4211
$ws.Range("L71").Value = 32169.7899
$ws.Range("N71").Value = -40281.7899
# Row 107
$ws.Range("H107").Value = 1935.1305
$ws.Range("I107").Value = 259.5
$ws.Range("J107").Value = 2094.7144
$ws.Range("K107").Value = 778.5
$ws.Range("L107").Value = 6284.1432
$ws.Range("M107").Value = 1141.5
$ws.Range("N107").Value = -10124.1432
# Row 121
$ws.Range("H121").Value = 83303.664
$ws.Range("I121").Value = 4791.6665
$ws.Range("J121").Value = 102931.664
$ws.Range("K121").Value = 14374.9995
$ws.Range("L121").Value = 308794.992
$ws.Range("M121").Value = -13064.9995
$ws.Range("N121").Value = -311414.992
# Row 122
$ws.Range("H122").Value = 1677.1538
$ws.Range("J122").Value = 1785.8182
$ws.Range("L122").Value = 16072.3638
$ws.Range("N122").Value = -20972.3638
# Row 130
$ws.Range("H130").Value = 2074.75
$ws.Range("I130").Value = 2074.75
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 6224.25
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -1204.25
$ws.Range("N130").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 139
$ws.Range("H139").Value = 66521.69500000001
$ws.Range("J139").Value = 66521.69500000001
$ws.Range("L139").Value = 66521.69500000001
$ws.Range("N139").Value = -76801.69500000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 3879.8333
$ws.Range("I107").Value = 4767.609
$ws.Range("J107").Value = 962.8570999999999
$ws.Range("K107").Value = 14302.827
$ws.Range("L107").Value = 2888.5713
$ws.Range("M107").Value = -12382.827
$ws.Range("N107").Value = -6728.5713
# Row 132
$ws.Range("H132").Value = 5252.533
$ws.Range("I132").Value = 4975.6665
$ws.Range("J132").Value = 6360
$ws.Range("K132").Value = 14926.9995
$ws.Range("L132").Value = 19080
$ws.Range("M132").Value = -12396.9995
$ws.Range("N132").Value = -24140

Write-Host "Cerberus_Profits market-data refresh applied."
